$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B20").Value = "Supports Stand-Up (20 pièces)"
$ws.Range("D20").Value = 50
$ws.Range("D20").NumberFormat = '#,##0.00\ "CHF"'

$ws.Range("D21").Select()
